$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row data: Sending cluster, Ligand symbol, Receptor symbol, Target cluster,
# then 16 numeric columns (E..T) recomputed to include the new "Neutro" sending
# cluster group (per Dr Hou's advice).
$rows = @(
    @("M1", "Npy", "Fap", "ECs", 3.0, 1.0, 2.709067, 8.127201, 0.670840434061904, 0.670840434061904, 3.0, 1.0, 1.262849333333333, 3.788548, 0.02753556474575134, 0.02753556474575133, 3.421143454905333, 30.790291094148, 0.01847197020617949, 0.01847197020617948),
    @("M1", "Npy", "Fap", "FAPs", 3.0, 1.0, 2.709067, 8.127201, 0.670840434061904, 0.670840434061904, 3.0, 1.0, 25.45149466666667, 76.354484, 0.5549524086300172, 0.554952408630017, 68.94980430214265, 620.548238719284, 0.3722845146890598, 0.3722845146890598),
    @("M1", "Npy", "Fap", "M1", 3.0, 1.0, 2.709067, 8.127201, 0.670840434061904, 0.670840434061904, 3.0, 1.0, 2.981677666666667, 8.945033, 0.06501343927129398, 0.06501343927129398, 8.077564571403666, 72.698081142633, 0.04361364382061209, 0.04361364382061209),
    @("M1", "Npy", "Fap", "M2", 3.0, 1.0, 2.709067, 8.127201, 0.670840434061904, 0.670840434061904, 3.0, 1.0, 2.257220666666667, 6.771662, 0.04921715059102959, 0.04921715059102958, 6.114962019784666, 55.034658178062, 0.03301685466577638, 0.03301685466577638),
    @("M1", "Npy", "Fap", "Neutro", 3.0, 1.0, 2.709067, 8.127201, 0.670840434061904, 0.670840434061904, 3.0, 1.0, 2.429399, 7.288196999999999, 0.05297138121868605, 0.05297138121868605, 6.581404660732998, 59.23264194659699, 0.03553534436960194, 0.03553534436960194),
    @("M1", "Npy", "Fap", "sCs", 3.0, 1.0, 2.709067, 8.127201, 0.670840434061904, 0.670840434061904, 3.0, 1.0, 11.47984033333333, 34.439521, 0.250310055543222, 0.250310055543222, 31.09965661230233, 279.896909510721, 0.1679181063106743, 0.1679181063106743),
    @("M2", "Npy", "Fap", "ECs", 3.0, 1.0, 1.267349666666667, 3.802049, 0.3138310719132735, 0.3138310719132735, 3.0, 1.0, 1.262849333333333, 3.788548, 0.02753556474575134, 0.02753556474575133, 1.600471681650223, 14.404245134852, 0.008641515799896488, 0.008641515799896485),
    @("M2", "Npy", "Fap", "FAPs", 3.0, 1.0, 1.267349666666667, 3.802049, 0.3138310719132735, 0.3138310719132735, 3.0, 1.0, 25.45149466666667, 76.354484, 0.5549524086300172, 0.554952408630017, 32.25594328196845, 290.303489537716, 0.1741613092612113, 0.1741613092612112),
    @("M2", "Npy", "Fap", "M1", 3.0, 1.0, 1.267349666666667, 3.802049, 0.3138310719132735, 0.3138310719132735, 3.0, 1.0, 2.981677666666667, 8.945033, 0.06501343927129398, 0.06501343927129398, 3.778828196957445, 34.009453772617, 0.0204032373352787, 0.0204032373352787),
    @("M2", "Npy", "Fap", "M2", 3.0, 1.0, 1.267349666666667, 3.802049, 0.3138310719132735, 0.3138310719132735, 3.0, 1.0, 2.257220666666667, 6.771662, 0.04921715059102959, 0.04921715059102958, 2.860687859493111, 25.746190735438, 0.01544587112649982, 0.01544587112649982),
    @("M2", "Npy", "Fap", "Neutro", 3.0, 1.0, 1.267349666666667, 3.802049, 0.3138310719132735, 0.3138310719132735, 3.0, 1.0, 2.429399, 7.288196999999999, 0.05297138121868605, 0.05297138121868605, 3.078898012850333, 27.710082115653, 0.01662406534858689, 0.01662406534858689),
    @("M2", "Npy", "Fap", "sCs", 3.0, 1.0, 1.267349666666667, 3.802049, 0.3138310719132735, 0.3138310719132735, 3.0, 1.0, 11.47984033333333, 34.439521, 0.250310055543222, 0.250310055543222, 14.54897181983656, 130.940746378529, 0.0785550730418004, 0.0785550730418004),
    @("Neutro", "Npy", "Fap", "ECs", 1.0, 0.3333333333333333, 0.06190133333333334, 0.185704, 0.01532849402482255, 0.01532849402482255, 3.0, 1.0, 1.262849333333333, 3.788548, 0.02753556474575134, 0.02753556474575133, 0.07817205753244445, 0.7035485177920001, 0.0004220787396753638, 0.0004220787396753637),
    @("Neutro", "Npy", "Fap", "FAPs", 1.0, 0.3333333333333333, 0.06190133333333334, 0.185704, 0.01532849402482255, 0.01532849402482255, 3.0, 1.0, 25.45149466666667, 76.354484, 0.5549524086300172, 0.554952408630017, 1.575481455192889, 14.179333096736, 0.008506584679746099, 0.008506584679746097),
    @("Neutro", "Npy", "Fap", "M1", 1.0, 0.3333333333333333, 0.06190133333333334, 0.185704, 0.01532849402482255, 0.01532849402482255, 3.0, 1.0, 2.981677666666667, 8.945033, 0.06501343927129398, 0.06501343927129398, 0.1845698231368889, 1.661128408232, 0.0009965581154031933, 0.0009965581154031933),
    @("Neutro", "Npy", "Fap", "M2", 1.0, 0.3333333333333333, 0.06190133333333334, 0.185704, 0.01532849402482255, 0.01532849402482255, 3.0, 1.0, 2.257220666666667, 6.771662, 0.04921715059102959, 0.04921715059102958, 0.1397249688942222, 1.257524720048, 0.0007544247987533886, 0.0007544247987533885),
    @("Neutro", "Npy", "Fap", "Neutro", 1.0, 0.3333333333333333, 0.06190133333333334, 0.185704, 0.01532849402482255, 0.01532849402482255, 3.0, 1.0, 2.429399, 7.288196999999999, 0.05297138121868605, 0.05297138121868605, 0.1503830372986666, 1.353447335688, 0.0008119715004972265, 0.0008119715004972264),
    @("Neutro", "Npy", "Fap", "sCs", 1.0, 0.3333333333333333, 0.06190133333333334, 0.185704, 0.01532849402482255, 0.01532849402482255, 3.0, 1.0, 11.47984033333333, 34.439521, 0.250310055543222, 0.250310055543222, 0.7106174230871112, 6.395556807784, 0.003836876190747279, 0.003836876190747279)
)

$startRow = 2
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $rowData = $rows[$i]
    for ($c = 0; $c -lt $rowData.Count; $c++) {
        $col = $c + 1
        $val = $rowData[$c]
        if ($col -le 4) {
            $ws.Cells.Item($r, $col).Value = $val
        } else {
            $ws.Cells.Item($r, $col).Value2 = $val
        }
    }
}